$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value of 45205 (2023-10-06) for
# all data rows (2 through 382). The edit bumps that date forward by one day
# to 45206 (2023-10-07) for every one of those rows.
$lastRow = 382
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45205) {
        $cell.Value2 = 45206
    }
}
